$p = $ppt.ActivePresentation
$s = $p.Slides.Item(5)
$shp = $s.Shapes.Item(2)
$tr = $shp.TextFrame.TextRange

# --- Change 1: merge the 3 runs that read
#     " for applications requires transparent federation of  capabilities " /
#     "and the interoperability " / "of services"
#     into a single run (text concatenated, keeping the first run's rPr). ---
$para1 = $tr.Paragraphs(3)
$merged = $para1.Characters(31, 103)
$merged.Text = " for applications requires transparent federation of  capabilities and the interoperability of services"

# --- Change 2: split the single run
#     "Test-beds exist out there but need to support effort to use test-beds
#      to transition from experimental to production-grade. "
#     into 4 runs, inserting "end-to-" so it reads "...support end-to-end
#     effort to use...". ---
$para2 = $tr.Paragraphs(8)

# Insert " end-to-end" right before the existing " effort" text (anchor on
# the space that precedes "effort" so the untouched tail keeps flowing from
# the original run).
$anchor = $para2.Characters(46, 1)
$anchor.Text = " end-to-end "

# Re-cut the run boundaries so the final shape is exactly:
#   "Test-beds exist out there but need to support" | " end-to" |
#   "-end effort " | "to use test-beds to transition from experimental to production-grade. "
$run2 = $para2.Characters(46, 7)
$run2.Text = $run2.Text

$run3 = $para2.Characters(53, 12)
$run3.Text = $run3.Text
